$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11:R11").Value = "highest potential for addressing the respective sustainable development goal"

$ws.Range("R11").Select()
